$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 33
$ws.Range("H33").Value = 273.92856
$ws.Range("I33").Value = 191.6
$ws.Range("K33").Value = 191.6
$ws.Range("M33").Value = 37.40000000000001

# ALC row 62
$ws.Range("H62").Value = 2348.5
$ws.Range("I62").Value = 2497
$ws.Range("J62").Value = 2200
$ws.Range("K62").Value = 2497
$ws.Range("L62").Value = 2200
$ws.Range("M62").Value = -1873
$ws.Range("N62").Value = -3448

# ALC row 65
$ws.Range("H65").Value = 2348.5
$ws.Range("I65").Value = 2497
$ws.Range("J65").Value = 2200
$ws.Range("K65").Value = 12485
$ws.Range("L65").Value = 11000
$ws.Range("M65").Value = -9365
$ws.Range("N65").Value = -17240

# ALC row 92
$ws.Range("H92").Value = 9500
$ws.Range("I92").Value = 9500
$ws.Range("K92").Value = 9500
$ws.Range("M92").Value = -8252

# ALC row 113
$ws.Range("H113").Value = 9044.546
$ws.Range("I113").Value = 8831.666999999999
$ws.Range("K113").Value = 8831.666999999999
$ws.Range("M113").Value = -5577.666999999999

# ALC row 116
$ws.Range("H116").Value = 3800
$ws.Range("I116").Value = 3800
$ws.Range("K116").Value = 3800
$ws.Range("M116").Value = -358

# ALC row 132
$ws.Range("H132").Value = 1728.7667
$ws.Range("I132").Value = 1456.3077
$ws.Range("K132").Value = 4368.9231
$ws.Range("M132").Value = -1838.9231

# ALC row 137
$ws.Range("H137").Value = 873.4
$ws.Range("I137").Value = 789
$ws.Range("K137").Value = 2367
$ws.Range("M137").Value = 183

# ALC row 138
$ws.Range("H138").Value = 3356.524
$ws.Range("I138").Value = 2884.7144
$ws.Range("J138").Value = 3592.4285
$ws.Range("K138").Value = 8654.143199999999
$ws.Range("L138").Value = 10777.2855
$ws.Range("M138").Value = -3514.143199999999
$ws.Range("N138").Value = -21057.2855

$ws = $wb.Worksheets.Item("ARM")
# ARM row 61
$ws.Range("H61").Value = 1489.5714
$ws.Range("I61").Value = 1489.5714
$ws.Range("K61").Value = 1489.5714
$ws.Range("M61").Value = -1277.5714

# ARM row 62
$ws.Range("H62").Value = 45000
$ws.Range("J62").Value = 45000
$ws.Range("L62").Value = 45000
$ws.Range("N62").Value = -46248

# ARM row 63
$ws.Range("H63").Value = 4583.1665

# ARM row 65
$ws.Range("H65").Value = 45000
$ws.Range("J65").Value = 45000
$ws.Range("L65").Value = 135000
$ws.Range("N65").Value = -141240

# ARM row 66
$ws.Range("H66").Value = 4583.1665

# ARM row 74
$ws.Range("H74").Value = 8000
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

# ARM row 77
$ws.Range("H77").Value = 8000
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

# ARM row 97
$ws.Range("H97").Value = 1516.0834
$ws.Range("I97").Value = 576
$ws.Range("J97").Value = 2832.2
$ws.Range("K97").Value = 576
$ws.Range("L97").Value = 2832.2
$ws.Range("M97").Value = -80
$ws.Range("N97").Value = -3824.2

# ARM row 136
$ws.Range("H136").Value = 1489.5714
$ws.Range("I136").Value = 1489.5714
$ws.Range("K136").Value = 4468.7142
$ws.Range("M136").Value = -1918.7142

$ws = $wb.Worksheets.Item("BSM")
# BSM row 20
$ws.Range("H20").Value = 2496.4
$ws.Range("I20").Value = 1157.6666
$ws.Range("K20").Value = 1157.6666
$ws.Range("M20").Value = -910.6666

# BSM row 105
$ws.Range("H105").Value = 3760.889
$ws.Range("I105").Value = 3293.5
$ws.Range("K105").Value = 3293.5
$ws.Range("M105").Value = -1546.5

# BSM row 107
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# CRP row 50
$ws.Range("H50").Value = 24444.445
$ws.Range("J50").Value = 24444.445
$ws.Range("L50").Value = 24444.445
$ws.Range("N50").Value = -25694.445

# CRP row 86
$ws.Range("H86").Value = 14979353
$ws.Range("I86").Value = 18399192
$ws.Range("K86").Value = 18399192
$ws.Range("M86").Value = -18398069

# CRP row 89
$ws.Range("H89").Value = 14979353
$ws.Range("I89").Value = 18399192
$ws.Range("K89").Value = 91995960
$ws.Range("M89").Value = -91990344

# CRP row 99
$ws.Range("H99").Value = 5287.75
$ws.Range("I99").Value = 5614.5713
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 5614.5713
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = -4116.5713
$ws.Range("N99").Value = -5996

# CRP row 105
$ws.Range("H105").Value = 17615.334
$ws.Range("I105").Value = 17615.334
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 17615.334
$ws.Range("L105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -15868.334

# CRP row 107
$ws.Range("H107").Value = 328.33334
$ws.Range("I107").Value = 297.8125
$ws.Range("J107").Value = 426
$ws.Range("K107").Value = 297.8125
$ws.Range("L107").Value = 426
$ws.Range("M107").Value = 1622.1875
$ws.Range("N107").Value = -4266

# CRP row 126
$ws.Range("H126").Value = 5287.75
$ws.Range("I126").Value = 5614.5713
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 16843.7139
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -14373.7139
$ws.Range("N126").Value = -13940

$ws = $wb.Worksheets.Item("CUL")
# CUL row 23
$ws.Range("H23").Value = 227.5
$ws.Range("I23").Value = 95
$ws.Range("K23").Value = 285
$ws.Range("M23").Value = -50

$ws = $wb.Worksheets.Item("GSM")
# GSM row 97
$ws.Range("H97").Value = 3333.3333
$ws.Range("I97").Value = 2750
$ws.Range("K97").Value = 2750
$ws.Range("M97").Value = -2254

# GSM row 102
$ws.Range("H102").Value = 2316.5
$ws.Range("I102").Value = 2185
$ws.Range("J102").Value = 3500
$ws.Range("K102").Value = 2185
$ws.Range("L102").Value = 3500
$ws.Range("M102").Value = -563
$ws.Range("N102").Value = -6744

$ws = $wb.Worksheets.Item("LTW")
# LTW row 16
$ws.Range("H16").Value = 612.875
$ws.Range("I16").Value = 691.8570999999999
$ws.Range("J16").Value = 60
$ws.Range("K16").Value = 691.8570999999999
$ws.Range("L16").Value = 60
$ws.Range("M16").Value = -521.8570999999999
$ws.Range("N16").Value = -400

$ws = $wb.Worksheets.Item("WVR")
# WVR row 62
$ws.Range("H62").Value = 7428.2856
$ws.Range("I62").Value = 7999.25
$ws.Range("J62").Value = 6667
$ws.Range("K62").Value = 7999.25
$ws.Range("L62").Value = 6667
$ws.Range("M62").Value = -7375.25
$ws.Range("N62").Value = -7915

# WVR row 65
$ws.Range("H65").Value = 7428.2856
$ws.Range("I65").Value = 7999.25
$ws.Range("J65").Value = 6667
$ws.Range("K65").Value = 39996.25
$ws.Range("L65").Value = 33335
$ws.Range("M65").Value = -36876.25
$ws.Range("N65").Value = -39575

# WVR row 81
$ws.Range("H81").Value = 3423
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

# WVR row 84
$ws.Range("H84").Value = 3423
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

# WVR row 135
$ws.Range("H135").Value = 45903.332
$ws.Range("J135").Value = 45903.332
$ws.Range("L135").Value = 45903.332
$ws.Range("N135").Value = -56043.332
